$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that currently sits right after
#    the "cliente" run, in the "Estimado/a <nombre de cliente>" paragraph.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) Split the big "Yo <nombre del cliente> confirmo que el servicio..."
#    run into three runs by inserting a new sentence in the middle:
#      "Yo <nombre del cliente> confirmo "
#      "que se me entrego <lista de software, hardware entregado> y "
#      "que el servicio fue <deficiente, excelente> ..."
#    and re-create the _GoBack bookmark (collapsed) between run 2 and
#    run 3, exactly where the cursor would have been left after typing
#    the inserted sentence. The two runs that trail the original
#    sentence ("y con ello confirmo..." / " apruebo que dicha carta...")
#    must stay untouched/unmerged, so their boundary is protected with
#    temporary bookmarks before the text is inserted.
# ---------------------------------------------------------------------

# Locate (dynamically, by content) the two run boundaries that already
# exist after the sentence we are splitting, so they survive the
# paragraph-level re-serialization triggered by the text insertion.
$protect1 = $d.Content
$protect1.Find.Execute("y con ello confirmo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$protectPos1 = $protect1.Start

$protect2 = $d.Content
$protect2.Find.Execute(" apruebo que dicha carta", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$protectPos2 = $protect2.Start

$d.Bookmarks.Add("_Protect1", $d.Range($protectPos1, $protectPos1))
$d.Bookmarks.Add("_Protect2", $d.Range($protectPos2, $protectPos2))

# Locate the insertion point: right before "que el servicio fue ...".
$anchorText = "que el servicio fue <deficiente, excelente> para la resolución de mi inconveniente presentado (<nombre de inconveniente>) "
$searchRange = $d.Content
$found = $searchRange.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target sentence to edit."
}

$insertPos = $searchRange.Start
$insertText = "que se me entrego <lista de software, hardware entregado> y "

$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertBefore($insertText)

$splitPoint1 = $insertPos
$splitPoint2 = $insertPos + $insertText.Length

# Temporary bookmark at the first split point forces the new text to
# live in its own run instead of being silently merged back into the
# preceding run; it is removed again once the split has taken effect.
$d.Bookmarks.Add("_TmpSplit", $d.Range($splitPoint1, $splitPoint1))

# The real bookmark goes at the second split point (between the newly
# inserted sentence and "que el servicio fue ...").
$d.Bookmarks.Add("_GoBack", $d.Range($splitPoint2, $splitPoint2))

$d.Bookmarks.Item("_TmpSplit").Delete()
$d.Bookmarks.Item("_Protect1").Delete()
$d.Bookmarks.Item("_Protect2").Delete()
